$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- New column L ("Maneuverability") in the "My version" table (rows 19-27) ---
$ws.Range("L19").Value = "Maneuverability"

$ws.Range("L20").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("L22").Value = 0.9
$ws.Range("L23").Value = 0.45
$ws.Range("L24").Value = 0.65
$ws.Range("L25").Value = 0.3
$ws.Range("L26").Value = 0.6
$ws.Range("L27").Value = 0.35

# --- Jolly boat row (20): Cannons columns changed ---
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0

# --- Column widths for G:L (best achievable precision through ColumnWidth) ---
$ws.Columns.Item(7).ColumnWidth = 12.333333333333334
$ws.Columns.Item(8).ColumnWidth = 5.0
$ws.Columns.Item(9).ColumnWidth = 14.333333333333334
$ws.Columns.Item(10).ColumnWidth = 12.333333333333334
$ws.Columns.Item(11).ColumnWidth = 11.5
$ws.Columns.Item(12).ColumnWidth = 15.0

# --- View: selection + scroll position ---
$ws.Activate()
$ws.Range("A10").Select()
$ws.Range("L18").Select()
